# Applies the "Finish 1 problem from 3.2" update to the time log workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D61 keeps its existing note text ("Finished 3 small problems, 15 mins on
# 5 problems from 3.1") - no change needed there.

# D62: update note to reflect the 1 problem finished from 3.2.
$ws.Range("D62").Value = "Finished 5 problems from 3.1, 1 problem from 3.2"

# C62: additional 0.25 hours logged (0.5 -> 0.75), which flows through the
# weekly SUM (D64) and the table SUBTOTAL (C65) formulas automatically.
$ws.Range("C62").Value = 0.75

# Update the active selection to match the saved workbook view (C63).
$ws.Range("C63").Select()
